$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 69.89967366666667
$ws.Range("H2").Value = 209.699021
$ws.Range("I2").Value = 0.6608367681537789
$ws.Range("J2").Value = 0.660836768153779
$ws.Range("M2").Value = 2.759544333333333
$ws.Range("N2").Value = 8.278632999999999
$ws.Range("O2").Value = 0.2574067337278401
$ws.Range("P2").Value = 0.2574067337278401
$ws.Range("Q2").Value = 192.8912483686992
$ws.Range("R2").Value = 1736.021235318293
$ws.Range("S2").Value = 0.1701038340177261
$ws.Range("T2").Value = 0.1701038340177262

$ws.Range("G3").Value = 69.89967366666667
$ws.Range("H3").Value = 209.699021
$ws.Range("I3").Value = 0.6608367681537789
$ws.Range("J3").Value = 0.660836768153779
$ws.Range("O3").Value = 0.6758254232987829
$ws.Range("P3").Value = 0.6758254232987829
$ws.Range("Q3").Value = 506.4390029409224
$ws.Range("R3").Value = 4557.951026468301
$ws.Range("S3").Value = 0.4466102885689273
$ws.Range("T3").Value = 0.4466102885689274

$ws.Range("G4").Value = 69.89967366666667
$ws.Range("H4").Value = 209.699021
$ws.Range("I4").Value = 0.6608367681537789
$ws.Range("J4").Value = 0.660836768153779
$ws.Range("M4").Value = 0.5200313333333334
$ws.Range("N4").Value = 1.560094
$ws.Range("O4").Value = 0.0485078515798926
$ws.Range("P4").Value = 0.0485078515798926
$ws.Range("Q4").Value = 36.35002049644157
$ws.Range("R4").Value = 327.150184467974
$ws.Range("S4").Value = 0.0320557718681394
$ws.Range("T4").Value = 0.0320557718681394

$ws.Range("G5").Value = 69.89967366666667
$ws.Range("H5").Value = 209.699021
$ws.Range("I5").Value = 0.6608367681537789
$ws.Range("J5").Value = 0.660836768153779
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1957573333333333
$ws.Range("N5").Value = 0.587272
$ws.Range("O5").Value = 0.01825999139348442
$ws.Range("P5").Value = 0.01825999139348442
$ws.Range("Q5").Value = 13.68337371785689
$ws.Range("R5").Value = 123.150363460712
$ws.Range("S5").Value = 0.01206687369898606
$ws.Range("T5").Value = 0.01206687369898607

$ws.Range("I6").Value = 0.1661491941864736
$ws.Range("J6").Value = 0.1661491941864736
$ws.Range("M6").Value = 2.759544333333333
$ws.Range("N6").Value = 8.278632999999999
$ws.Range("O6").Value = 0.2574067337278401
$ws.Range("P6").Value = 0.2574067337278401
$ws.Range("Q6").Value = 48.49718875603555
$ws.Range("R6").Value = 436.4746988043199
$ws.Range("S6").Value = 0.04276792138705281
$ws.Range("T6").Value = 0.04276792138705281

$ws.Range("I7").Value = 0.1661491941864736
$ws.Range("J7").Value = 0.1661491941864736
$ws.Range("O7").Value = 0.6758254232987829
$ws.Range("P7").Value = 0.6758254232987829
$ws.Range("S7").Value = 0.1122878494918252
$ws.Range("T7").Value = 0.1122878494918252

$ws.Range("I8").Value = 0.1661491941864736
$ws.Range("J8").Value = 0.1661491941864736
$ws.Range("M8").Value = 0.5200313333333334
$ws.Range("N8").Value = 1.560094
$ws.Range("O8").Value = 0.0485078515798926
$ws.Range("P8").Value = 0.0485078515798926
$ws.Range("Q8").Value = 9.13921092952889
$ws.Range("R8").Value = 82.25289836576
$ws.Range("S8").Value = 0.008059540451716217
$ws.Range("T8").Value = 0.008059540451716217

$ws.Range("I9").Value = 0.1661491941864736
$ws.Range("J9").Value = 0.1661491941864736
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.1957573333333333
$ws.Range("N9").Value = 0.587272
$ws.Range("O9").Value = 0.01825999139348442
$ws.Range("P9").Value = 0.01825999139348442
$ws.Range("Q9").Value = 3.440307238542222
$ws.Range("R9").Value = 30.96276514688
$ws.Range("S9").Value = 0.00303388285587938
$ws.Range("T9").Value = 0.003033882855879381

$ws.Range("G10").Value = 4.152730666666667
$ws.Range("H10").Value = 12.458192
$ws.Range("I10").Value = 0.0392602278210887
$ws.Range("J10").Value = 0.03926022782108871
$ws.Range("M10").Value = 2.759544333333333
$ws.Range("N10").Value = 8.278632999999999
$ws.Range("O10").Value = 0.2574067337278401
$ws.Range("P10").Value = 0.2574067337278401
$ws.Range("Q10").Value = 11.45964437905955
$ws.Range("R10").Value = 103.136799411536
$ws.Range("S10").Value = 0.01010584700883732
$ws.Range("T10").Value = 0.01010584700883732

$ws.Range("G11").Value = 4.152730666666667
$ws.Range("H11").Value = 12.458192
$ws.Range("I11").Value = 0.0392602278210887
$ws.Range("J11").Value = 0.03926022782108871
$ws.Range("O11").Value = 0.6758254232987829
$ws.Range("P11").Value = 0.6758254232987829
$ws.Range("Q11").Value = 30.08747634986133
$ws.Range("R11").Value = 270.787287148752
$ws.Range("S11").Value = 0.02653306008599392
$ws.Range("T11").Value = 0.02653306008599393

$ws.Range("G12").Value = 4.152730666666667
$ws.Range("H12").Value = 12.458192
$ws.Range("I12").Value = 0.0392602278210887
$ws.Range("J12").Value = 0.03926022782108871
$ws.Range("M12").Value = 0.5200313333333334
$ws.Range("N12").Value = 1.560094
$ws.Range("O12").Value = 0.0485078515798926
$ws.Range("P12").Value = 0.0485078515798926
$ws.Range("Q12").Value = 2.159550065560889
$ws.Range("R12").Value = 19.435950590048
$ws.Range("S12").Value = 0.001904429304138141
$ws.Range("T12").Value = 0.001904429304138141

$ws.Range("G13").Value = 4.152730666666667
$ws.Range("H13").Value = 12.458192
$ws.Range("I13").Value = 0.0392602278210887
$ws.Range("J13").Value = 0.03926022782108871
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.1957573333333333
$ws.Range("N13").Value = 0.587272
$ws.Range("O13").Value = 0.01825999139348442
$ws.Range("P13").Value = 0.01825999139348442
$ws.Range("Q13").Value = 0.8129274813582222
$ws.Range("R13").Value = 7.316347332224001
$ws.Range("S13").Value = 0.0007168914221193172
$ws.Range("T13").Value = 0.0007168914221193175

$ws.Range("G14").Value = 14.14774133333333
$ws.Range("H14").Value = 42.443224
$ws.Range("I14").Value = 0.1337538098386587
$ws.Range("J14").Value = 0.1337538098386588
$ws.Range("M14").Value = 2.759544333333333
$ws.Range("N14").Value = 8.278632999999999
$ws.Range("O14").Value = 0.2574067337278401
$ws.Range("P14").Value = 0.2574067337278401
$ws.Range("Q14").Value = 39.04131942586577
$ws.Range("R14").Value = 351.3718748327919
$ws.Range("S14").Value = 0.03442913131422378
$ws.Range("T14").Value = 0.03442913131422379

$ws.Range("G15").Value = 14.14774133333333
$ws.Range("H15").Value = 42.443224
$ws.Range("I15").Value = 0.1337538098386587
$ws.Range("J15").Value = 0.1337538098386588
$ws.Range("O15").Value = 0.6758254232987829
$ws.Range("P15").Value = 0.6758254232987829
$ws.Range("Q15").Value = 102.5035974972827
$ws.Range("R15").Value = 922.532377475544
$ws.Range("S15").Value = 0.09039422515203645
$ws.Range("T15").Value = 0.09039422515203648

$ws.Range("G16").Value = 14.14774133333333
$ws.Range("H16").Value = 42.443224
$ws.Range("I16").Value = 0.1337538098386587
$ws.Range("J16").Value = 0.1337538098386588
$ws.Range("M16").Value = 0.5200313333333334
$ws.Range("N16").Value = 1.560094
$ws.Range("O16").Value = 0.0485078515798926
$ws.Range("P16").Value = 0.0485078515798926
$ws.Range("Q16").Value = 7.357268789228446
$ws.Range("R16").Value = 66.21541910305601
$ws.Range("S16").Value = 0.006488109955898836
$ws.Range("T16").Value = 0.006488109955898837

$ws.Range("G17").Value = 14.14774133333333
$ws.Range("H17").Value = 42.443224
$ws.Range("I17").Value = 0.1337538098386587
$ws.Range("J17").Value = 0.1337538098386588
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.1957573333333333
$ws.Range("N17").Value = 0.587272
$ws.Range("O17").Value = 0.01825999139348442
$ws.Range("P17").Value = 0.01825999139348442
$ws.Range("Q17").Value = 2.769524116103111
$ws.Range("R17").Value = 24.925717044928
$ws.Range("S17").Value = 0.00244234341649966
$ws.Range("T17").Value = 0.002442343416499661
